$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top; this pushes all existing rows (1-24) down to (2-25)
$ws.Rows.Item(1).Insert()

# Copy formatting (style) from the (now shifted) header cells A2:L2 onto the new A1:L1
$ws.Range("A2:L2").Copy()
$ws.Range("A1:L1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate new row 1 with a sequential numeric index header (0-based column index)
$ws.Range("A1").Value = 0
$ws.Range("B1").Value = 1
$ws.Range("C1").Value = 2
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 4
$ws.Range("F1").Value = 5
$ws.Range("G1").Value = 6
$ws.Range("H1").Value = 7
$ws.Range("I1").Value = 8
$ws.Range("J1").Value = 9
$ws.Range("K1").Value = 10
$ws.Range("L1").Value = 11
